$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.321.94"
$ws.Range("E2").Value = "  -4.34%  "

$ws.Range("D3").Value = "1.859.60"
$ws.Range("E3").Value = "  -5.34%  "

$ws.Range("E4").Value = "  -0.98%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.21%  "

$ws.Range("E6").Value = "  -0.88%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4509"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -6.73%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3860"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -5.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.11"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -10.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07897"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.020"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.39"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.08%  "

$ws.Range("D13").Value = "1.851.86"
$ws.Range("E13").Value = "  -9.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.161"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.16%  "

$ws.Range("E15").Value = "  -5.41%  "

$ws.Range("E16").Value = "  -1.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001034"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "85.62"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -6.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06531"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.04"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -8.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.515"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -6.30%  "

$ws.Range("D23").Value = "27.336.71"
$ws.Range("E23").Value = "  -4.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.78"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -6.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.268"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.30%  "

$ws.Range("D26").Value = "2.090.00"
$ws.Range("E26").Value = "  -8.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.96"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.73"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.064"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.494"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -7.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.63"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.484"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09292"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9369"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -6.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.599"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.275"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02230"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05991"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.208"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.256"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -10.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5900"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.14"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -10.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.259"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5633"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.84%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.98"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -9.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.922"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -7.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.360"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.59%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06798"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "108.22"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.07%  "

